# Add Sunday's round of scores to the "jul18" sheet and make it the active tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("jul18")

# New round date (row 22) - Sunday 7/7/2018, formatted like the other date cells.
$ws.Range("A22").Value = 43288
$ws.Range("A22").NumberFormat = $ws.Range("A1").NumberFormat

# Hole-by-hole scores for the round played on Sunday.
$scores = @(
    @{ Row = 23; Hole = "Hole 1";  Score = 5; Note = "S"; Putts = 3 },
    @{ Row = 24; Hole = "Hole 2";  Score = 4; Note = $null; Putts = 2 },
    @{ Row = 25; Hole = "Hole 3";  Score = 4; Note = "R"; Putts = 2 },
    @{ Row = 26; Hole = "Hole 4";  Score = 5; Note = "S"; Putts = 2 },
    @{ Row = 27; Hole = "Hole 5";  Score = 4; Note = $null; Putts = 3 },
    @{ Row = 28; Hole = "Hole 6";  Score = 6; Note = "R"; Putts = 3 },
    @{ Row = 29; Hole = "Hole 7";  Score = 5; Note = "R"; Putts = 2 },
    @{ Row = 30; Hole = "Hole 8";  Score = 4; Note = $null; Putts = 1 },
    @{ Row = 31; Hole = "Hole 9";  Score = 5; Note = "S"; Putts = 2 },
    @{ Row = 32; Hole = "Hole 10"; Score = 4; Note = "S"; Putts = 2 },
    @{ Row = 33; Hole = "Hole 11"; Score = 4; Note = $null; Putts = 2 },
    @{ Row = 34; Hole = "Hole 12"; Score = 5; Note = "R"; Putts = 2 },
    @{ Row = 35; Hole = "Hole 13"; Score = 6; Note = "S"; Putts = 3 },
    @{ Row = 36; Hole = "Hole 14"; Score = 5; Note = "S"; Putts = 2 },
    @{ Row = 37; Hole = "Hole 15"; Score = 4; Note = $null; Putts = 3 },
    @{ Row = 38; Hole = "Hole 16"; Score = 5; Note = "R"; Putts = 2 },
    @{ Row = 39; Hole = "Hole 17"; Score = 4; Note = "S"; Putts = 2 },
    @{ Row = 40; Hole = "Hole 18"; Score = 5; Note = "R"; Putts = 2 }
)

foreach ($entry in $scores) {
    $r = $entry.Row
    $ws.Range("A$r").Value = $entry.Hole
    $ws.Range("B$r").Value = $entry.Score
    if ($entry.Note) {
        $ws.Range("C$r").Value = $entry.Note
    }
    $ws.Range("E$r").Value = $entry.Putts
}

# Totals row for the new round.
$ws.Range("B41").Formula = "=SUM(B23:B40)"
$ws.Range("E41").Formula = "=SUM(E23:E40)"

# Make "jul18" the active sheet/tab and put the selection on the new Putts entry,
# matching where focus was left after entering Sunday's scores.
$ws.Activate() | Out-Null
$ws.Range("E23").Select() | Out-Null
